# WRESBAL.xlsx refresh: roll the weekly FRED series forward.
# - Drop the 61 oldest observations (old rows 2-62, dates 44048..44468)
#   so the series window starts at 44475 (2021-10-06) instead of 44048.
# - Append 22 new weekly observations through 45112 (2023-07-05).
# - Refresh the SeriesInfo metadata sheet (realtime_start/end,
#   observation_end, last_updated, popularity).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Data" sheet: drop the oldest 61 rows, then append new rows.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Remove the oldest 61 observations (old rows 2..62); this shifts the
# remaining rows (old 63..132) up to become the new rows 2..71.
$ws.Range("A2:A62").EntireRow.Delete() | Out-Null

# New weekly observations to append after the shifted data (new rows
# 72..93, continuing the series from 44965 through 45112).
$newRows = @(
    @{Row=72; Date=44965; Value=3037.584},
    @{Row=73; Date=44972; Value=3041.038},
    @{Row=74; Date=44979; Value=3016.577},
    @{Row=75; Date=44986; Value=2998.303},
    @{Row=76; Date=44993; Value=2999.72},
    @{Row=77; Date=45000; Value=3251.482},
    @{Row=78; Date=45007; Value=3425.075},
    @{Row=79; Date=45014; Value=3437.592},
    @{Row=80; Date=45021; Value=3292.058},
    @{Row=81; Date=45028; Value=3413.989},
    @{Row=82; Date=45035; Value=3303.152},
    @{Row=83; Date=45042; Value=3135.655},
    @{Row=84; Date=45049; Value=3089.706},
    @{Row=85; Date=45056; Value=3196.429},
    @{Row=86; Date=45063; Value=3248.878},
    @{Row=87; Date=45070; Value=3240.172},
    @{Row=88; Date=45077; Value=3309.447},
    @{Row=89; Date=45084; Value=3349.751},
    @{Row=90; Date=45091; Value=3325.261},
    @{Row=91; Date=45098; Value=3251.558},
    @{Row=92; Date=45105; Value=3171.567},
    @{Row=93; Date=45112; Value=3114.989}
)

# Carry the date column's formatting (border/bold/center + the
# YYYY-MM-DD HH:MM:SS number format) down onto the freshly appended
# rows, the same way dragging/filling the column would in the UI.
$lastRow = $newRows[0].Row - 1
$ws.Range("A" + $lastRow).Copy() | Out-Null
$ws.Range("A" + $newRows[0].Row + ":A" + $newRows[-1].Row).PasteSpecial(-4122) | Out-Null

foreach ($row in $newRows) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Date
    $ws.Cells.Item($row.Row, 2).Value = $row.Value
}

# ---------------------------------------------------------------
# 2) "SeriesInfo" sheet: refresh the FRED metadata fields.
# ---------------------------------------------------------------
$info = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end / observation_end hold plain
# "YYYY-MM-DD" text (not real dates) in this sheet. Assigning such a
# string straight to .Value gets auto-coerced into a date serial by
# COM's smart-typing, so round-trip it through a quoted formula and
# paste-values-back to land a literal text value instead (keeps the
# default, un-styled cell formatting too).
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

Set-TextValue $info.Range("B3") "2023-07-09"
Set-TextValue $info.Range("B4") "2023-07-09"
Set-TextValue $info.Range("B7") "2023-07-05"

# last_updated already includes a time + UTC-offset suffix, so it is
# never mistaken for a date/number by COM - plain assignment is fine.
$info.Range("B14").Value = "2023-07-06 15:35:25-05"
$info.Range("B15").Value = 78
